$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh prepends a new week's worth of rows (2 rows: "1a
# (cosecha)" / "2a (cosecha)") to the top of the data block, pushing every
# existing data row down by two. Insert two blank rows right before the
# first data row (row 415; row 414 is the last row of the previous block
# end and stays put).
$ws.Range("A415:A416").EntireRow.Insert()

# New row 415: Cebolla, 1a (cosecha), Región de O'Higgins, fecha 2021-12-13
$ws.Cells.Item(415, 1).Value = 8
$ws.Cells.Item(415, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(415, 3).Value = "Coquimbo"
$ws.Cells.Item(415, 4).Value = 44543
$ws.Cells.Item(415, 5).Value = 4
$ws.Cells.Item(415, 6).Value = 100112004
$ws.Cells.Item(415, 7).Value = "Cebolla"
$ws.Cells.Item(415, 8).Value = "Sin especificar"
$ws.Cells.Item(415, 9).Value = "1a (cosecha)"
$ws.Cells.Item(415, 10).Value = 3200
$ws.Cells.Item(415, 11).Value = 4500
$ws.Cells.Item(415, 12).Value = 5000
$ws.Cells.Item(415, 13).Value = 4750
$ws.Cells.Item(415, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(415, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(415, 16).Value = 264
$ws.Cells.Item(415, 17).Value = 18
$ws.Cells.Item(415, 18).Value = "Hortaliza"

# New row 416: Cebolla, 2a (cosecha), Región de O'Higgins, fecha 2021-12-13
$ws.Cells.Item(416, 1).Value = 8
$ws.Cells.Item(416, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(416, 3).Value = "Coquimbo"
$ws.Cells.Item(416, 4).Value = 44543
$ws.Cells.Item(416, 5).Value = 4
$ws.Cells.Item(416, 6).Value = 100112004
$ws.Cells.Item(416, 7).Value = "Cebolla"
$ws.Cells.Item(416, 8).Value = "Sin especificar"
$ws.Cells.Item(416, 9).Value = "2a (cosecha)"
$ws.Cells.Item(416, 10).Value = 1600
$ws.Cells.Item(416, 11).Value = 4000
$ws.Cells.Item(416, 12).Value = 4200
$ws.Cells.Item(416, 13).Value = 4100
$ws.Cells.Item(416, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(416, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(416, 16).Value = 228
$ws.Cells.Item(416, 17).Value = 18
$ws.Cells.Item(416, 18).Value = "Hortaliza"
